$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update cell texts (use case steps renumbered / reworded) ---
# The order in which NEW (previously unused) strings are assigned matters,
# since each new unique string value is appended to the shared-strings table
# in assignment order. This order was derived to reproduce the target
# shared-strings layout.

# Row 6 / C6: step 1 description reworded
$ws.Cells.Item(6, 3).Value = "1. Apresenta as opções de seleção dos ingredientes e temperatura preferidos"

# Row 7 / B7: step 2 description reworded
$ws.Cells.Item(7, 2).Value = "2. Seleciona ingredientes/temperatura"

# Row 12 / A12: Alternative 1 heading, now references Passo 2 instead of Passo 3
$ws.Cells.Item(12, 1).Value = "Alternativa 1 [Não existem receitas para a combinação definida] (Passo 2)"

# Row 12 / C12: sub-step renumbered from 3.1 to 2.1
$ws.Cells.Item(12, 3).Value = "2.1 Informa que não existem receitas para a configuração"

# Row 13 / C13: sub-step renumbered from 3.2 to 2.2
$ws.Cells.Item(13, 3).Value = "2.2 Regressa ao Passo 1"

# Row 9 / B9: step 4 description reworded
$ws.Cells.Item(9, 2).Value = "4. Termina a Configuração Inicial  "

# Row 16 / C16: sub-step 4.2 now returns to step 1 instead of step 3
$ws.Cells.Item(16, 3).Value = "4.2 Regressa ao passo 1"

# --- Row height: row 6 shrinks now that its text is a single shorter line ---
$ws.Rows.Item(6).RowHeight = 38.25

# --- Sheet view: zoom to 85% and move selection to C17 ---
$ws.Activate() | Out-Null
$ws.Range("C17").Select() | Out-Null
$win = $excel.ActiveWindow
$win.Zoom = 85
